$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.903.86'
$ws.Range('D3').Value = '''3.526.41'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''598.67'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').Value = '''143.74'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').Value = '''3.523.75'
$ws.Range('E7').Value = '  -1.01%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').Value = '''0.497'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').Value = '''4.126.98'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').Value = '''0.0000200'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').Value = '''28.71'
$ws.Range('E15').Value = '  -4.55%  '
$ws.Range('D16').Value = '''3.521.01'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '''65.895.16'
$ws.Range('D19').Value = '''10.90'
$ws.Range('E19').Value = '  -5.55%  '
$ws.Range('D20').Value = '''6.18'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D21').Value = '''14.34'
$ws.Range('E21').Value = '  -3.97%  '
$ws.Range('D22').Value = '''414.50'
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('D24').Value = '''77.32'
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('D25').Value = '''3.669.83'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -3.04%  '
$ws.Range('D28').Value = '''7.79'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '''2.44'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''8.93'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '''3.525.07'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').Value = '''24.37'
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('D36').Value = '''7.52'
$ws.Range('E36').Value = '  -4.36%  '
$ws.Range('D37').Value = '''1.29'
$ws.Range('E37').Value = '  -12.12%  '
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('D39').Value = '''5.26'
$ws.Range('E39').Value = '  -6.38%  '
$ws.Range('E40').Value = '  -8.18%  '
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').Value = '''5.06'
$ws.Range('E42').Value = '  -2.65%  '
$ws.Range('D43').Value = '''0.858'
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('D44').Value = '''45.29'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('D45').Value = '''1.79'
$ws.Range('E45').Value = '  -7.78%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = '''2.41'
$ws.Range('E47').Value = '  -4.55%  '
$ws.Range('D48').Value = '''7.07'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  -7.04%  '
$ws.Range('D50').Value = '''22.56'
$ws.Range('D51').Value = '''23.06'
$ws.Range('E51').Value = '  -8.27%  '
